$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '21.334.18'
$ws.Range('E2').Value = '  +4.27%  '

$ws.Range('D3').Value = '1.549.21'
$ws.Range('E3').Value = '  +5.34%  '

$ws.Range('D4').Value = '1.007'
$ws.Range('E4').Value = '  -0.18%  '

$ws.Range('D5').Value = '0.9690'
$ws.Range('E5').Value = '  -0.01%  '

$ws.Range('D6').Value = '283.09'
$ws.Range('E6').Value = '  +2.70%  '

$ws.Range('D7').Value = '0.3618'
$ws.Range('E7').Value = '  -0.83%  '

$ws.Range('D8').Value = '0.3210'
$ws.Range('E8').Value = '  +4.71%  '

$ws.Range('D9').Value = '40.96'
$ws.Range('E9').Value = '  +2.71%  '

$ws.Range('D10').Value = '1.111'
$ws.Range('E10').Value = '  +6.15%  '

$ws.Range('D11').Value = '0.06939'
$ws.Range('E11').Value = '  +5.00%  '

$ws.Range('D12').Value = '1.002'

$ws.Range('D13').Value = '5.731'
$ws.Range('E13').Value = '  +5.25%  '

$ws.Range('D14').Value = '18.95'
$ws.Range('E14').Value = '  +4.99%  '

$ws.Range('D15').Value = '6.432'
$ws.Range('E15').Value = '  +4.56%  '

$ws.Range('D16').Value = '0.00001054'
$ws.Range('E16').Value = '  +2.47%  '

$ws.Range('D17').Value = '0.9687'
$ws.Range('E17').Value = '  -1.06%  '

$ws.Range('D18').Value = '1.544.27'
$ws.Range('E18').Value = '  +4.84%  '

$ws.Range('D19').Value = '0.06157'
$ws.Range('E19').Value = '  +4.55%  '

$ws.Range('D20').Value = '73.17'
$ws.Range('E20').Value = '  +5.97%  '

$ws.Range('D21').Value = '5.750'
$ws.Range('E21').Value = '  +5.57%  '

$ws.Range('E22').Value = '  +6.31%  '

$ws.Range('D23').Value = '11.42'
$ws.Range('E23').Value = '  +4.50%  '

$ws.Range('D24').Value = '2.323'
$ws.Range('E24').Value = '  +3.28%  '

$ws.Range('D25').Value = '21.344.37'
$ws.Range('E25').Value = '  +4.04%  '

$ws.Range('D26').Value = '147.80'
$ws.Range('E26').Value = '  +4.21%  '

$ws.Range('D27').Value = '2.277'
$ws.Range('E27').Value = '  +6.62%  '

$ws.Range('D28').Value = '17.89'
$ws.Range('E28').Value = '  +4.01%  '

$ws.Range('D29').Value = '1.714.64'
$ws.Range('E29').Value = '  +5.25%  '

$ws.Range('D30').Value = '118.68'
$ws.Range('E30').Value = '  +4.53%  '

$ws.Range('D31').Value = '4.050'
$ws.Range('E31').Value = '  +4.34%  '

$ws.Range('D32').Value = '0.8736'
$ws.Range('E32').Value = '  +9.18%  '

$ws.Range('D33').Value = '5.276'
$ws.Range('E33').Value = '  +6.66%  '

$ws.Range('D34').Value = '0.08058'
$ws.Range('E34').Value = '  +2.44%  '

$ws.Range('D35').Value = '1.534'
$ws.Range('E35').Value = '  +0.56%  '

$ws.Range('D36').Value = '4.999'
$ws.Range('E36').Value = '  +5.40%  '

$ws.Range('D37').Value = '1.212'
$ws.Range('E37').Value = '  -0.25%  '

$ws.Range('D38').Value = '0.05870'
$ws.Range('E38').Value = '  +2.48%  '

$ws.Range('B39').Value = 'Algorand'
$ws.Range('C39').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D39').Value = '0.2008'
$ws.Range('E39').Value = '  +6.93%  '

$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').Value = '0.02128'
$ws.Range('E40').Value = '  +4.63%  '

$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D41').Value = '7.966'
$ws.Range('E41').Value = '  +4.43%  '

$ws.Range('B42').Value = 'Aptos'
$ws.Range('C42').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D42').Value = '10.79'
$ws.Range('E42').Value = '  +3.45%  '

$ws.Range('D43').Value = '0.9679'
$ws.Range('E43').Value = '  -0.53%  '

$ws.Range('D44').Value = '0.5515'
$ws.Range('E44').Value = '  +4.46%  '

$ws.Range('D45').Value = '12.59'
$ws.Range('E45').Value = '  +4.32%  '

$ws.Range('D46').Value = '3.579'
$ws.Range('E46').Value = '  +2.38%  '

$ws.Range('D47').Value = '0.5516'
$ws.Range('E47').Value = '  +6.77%  '

$ws.Range('D48').Value = '122.28'
$ws.Range('E48').Value = '  +4.41%  '

$ws.Range('D49').Value = '1.883'
$ws.Range('E49').Value = '  +6.65%  '

$ws.Range('D50').Value = '0.06625'
$ws.Range('E50').Value = '  +2.80%  '

$ws.Range('D51').Value = '70.22'
$ws.Range('E51').Value = '  +4.99%  '

